$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the data mistake in "Fraud and Business Process Analytics" (row 76) ---
# Column G = Last_7_Weeks, was incorrectly marked as 1 (true), corrected to 0 (false)
$ws.Range("G76").Value = 0
# Undo the incidental row-height stamp that a value write triggers, so the row
# keeps using the sheet's default (non-custom) height, as in the original file.
$ws.Rows(76).AutoFit()

# --- Hide the rows that no longer match the active filter criteria ---
$rowsToHide = @(4, 12, 18, 19, 27, 28, 38, 39, 40, 50, 51, 55, 56, 57, 64, 71, 89, 90, 99)
foreach ($r in $rowsToHide) {
    $ws.Rows($r).Hidden = $true
}

# Row 92 ("Production Control") now matches the new filter, so make it visible
$ws.Rows(92).Hidden = $false

# --- Update the AutoFilter: filter on Class (column C) = "Production Control" ---
# instead of the previous filter on Start (column H) = 12.50
$ws.AutoFilterMode = $false
$ws.Range("A1:M99").AutoFilter(3, @("Production Control"), 7)

# --- Update the remembered selection (active cell moved from G101 to G100) ---
$ws.Range("G100").Select()
